$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.737.82"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "1.878.28"

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "331.99"
$ws.Range("E5").Value = "  +3.48%  "

# Row 6
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  +5.48%  "

# Row 8
$ws.Range("D8").Value = "0.3963"
$ws.Range("E8").Value = "  +3.29%  "

# Row 9
$ws.Range("D9").Value = "47.87"
$ws.Range("E9").Value = "  -0.18%  "

# Row 10
$ws.Range("D10").Value = "0.08037"
$ws.Range("E10").Value = "  +2.77%  "

# Row 11
$ws.Range("D11").Value = "1.024"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").Value = "21.85"
$ws.Range("E12").Value = "  +2.24%  "

# Row 13
$ws.Range("D13").Value = "1.868.18"
$ws.Range("E13").Value = "  +1.42%  "

# Row 14
$ws.Range("D14").Value = "5.966"
$ws.Range("E14").Value = "  +2.09%  "

# Row 15
$ws.Range("D15").Value = "7.157"
$ws.Range("E15").Value = "  +0.88%  "

# Row 16
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.15%  "

# Row 17
$ws.Range("D17").Value = "87.16"
$ws.Range("E17").Value = "  +1.82%  "

# Row 18
$ws.Range("D18").Value = "0.00001048"
$ws.Range("E18").Value = "  +2.50%  "

# Row 19
$ws.Range("D19").Value = "0.06612"
$ws.Range("E19").Value = "  +1.83%  "

# Row 20
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  +1.20%  "

# Row 21
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").Value = "27.692.17"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23
$ws.Range("D23").Value = "5.509"
$ws.Range("E23").Value = "  +0.69%  "

# Row 24
$ws.Range("E24").Value = "  +2.54%  "

# Row 25
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +0.56%  "

# Row 26
$ws.Range("D26").Value = "2.095.58"
$ws.Range("E26").Value = "  +4.33%  "

# Row 27
$ws.Range("D27").Value = "156.62"
$ws.Range("E27").Value = "  +3.81%  "

# Row 28
$ws.Range("D28").Value = "20.27"
$ws.Range("E28").Value = "  +4.83%  "

# Row 29
$ws.Range("D29").Value = "2.094"
$ws.Range("E29").Value = "  +3.43%  "

# Row 30
$ws.Range("D30").Value = "5.584"
$ws.Range("E30").Value = "  +2.27%  "

# Row 31
$ws.Range("D31").Value = "122.54"
$ws.Range("E31").Value = "  +2.45%  "

# Row 32
$ws.Range("D32").Value = "0.9681"
$ws.Range("E32").Value = "  +4.57%  "

# Row 33
$ws.Range("D33").Value = "0.09555"
$ws.Range("E33").Value = "  +2.35%  "

# Row 34
$ws.Range("E34").Value = "  -2.87%  "

# Row 35
$ws.Range("D35").Value = "3.626"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$ws.Range("D36").Value = "5.297"
$ws.Range("E36").Value = "  +1.67%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02267"
$ws.Range("E37").Value = "  +2.33%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06115"
$ws.Range("E38").Value = "  +3.05%  "

# Row 39
$ws.Range("D39").Value = "1.231"
$ws.Range("E39").Value = "  +0.65%  "

# Row 40
$ws.Range("D40").Value = "8.191"
$ws.Range("E40").Value = "  -0.97%  "

# Row 41
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5989"
$ws.Range("E42").Value = "  +1.87%  "

# Row 43
$ws.Range("D43").Value = "0.1902"
$ws.Range("E43").Value = "  +2.97%  "

# Row 44
$ws.Range("D44").Value = "10.24"
$ws.Range("E44").Value = "  +0.07%  "

# Row 45
$ws.Range("D45").Value = "1.261"
$ws.Range("E45").Value = "  -0.96%  "

# Row 46
$ws.Range("D46").Value = "0.5688"
$ws.Range("E46").Value = "  +1.61%  "

# Row 47
$ws.Range("D47").Value = "12.28"
$ws.Range("E47").Value = "  +0.95%  "

# Row 48
$ws.Range("D48").Value = "3.403"
$ws.Range("E48").Value = "  +1.51%  "

# Row 49
$ws.Range("E49").Value = "  +1.50%  "

# Row 50
$ws.Range("E50").Value = "  -0.25%  "

# Row 51
$ws.Range("E51").Value = "  +9.19%  "
